$wb = $excel.ActiveWorkbook

# --- Sheet: Summary ---
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B2").Value = 0.6282771535580525
$wsSummary.Range("C2").Value = 0.5828295042321644
$wsSummary.Range("D2").Value = 0.9026217228464419
$wsSummary.Range("E2").Value = 0.7083027185892726
$wsSummary.Range("F2").Value = 0.8133648329395883
$wsSummary.Range("G2").Value = 0.8839669887846512
$wsSummary.Range("H2").Value = 0.7739044593134985
$wsSummary.Range("I2").Value = 482
$wsSummary.Range("J2").Value = 345
$wsSummary.Range("K2").Value = 189
$wsSummary.Range("L2").Value = 52

# --- Sheet: Classification Report ---
$wsClass = $wb.Worksheets.Item("Classification Report")
$wsClass.Range("B2").Value = 0.7842323651452282
$wsClass.Range("C2").Value = 0.3539325842696629
$wsClass.Range("D2").Value = 0.487741935483871

$wsClass.Range("B3").Value = 0.5828295042321644
$wsClass.Range("C3").Value = 0.9026217228464419
$wsClass.Range("D3").Value = 0.7083027185892726

$wsClass.Range("B4").Value = 0.6282771535580525
$wsClass.Range("C4").Value = 0.6282771535580525
$wsClass.Range("D4").Value = 0.6282771535580525
$wsClass.Range("E4").Value = 0.6282771535580525

$wsClass.Range("B5").Value = 0.6835309346886963
$wsClass.Range("C5").Value = 0.6282771535580525
$wsClass.Range("D5").Value = 0.5980223270365718

$wsClass.Range("B6").Value = 0.6835309346886964
$wsClass.Range("C6").Value = 0.6282771535580525
$wsClass.Range("D6").Value = 0.5980223270365718

# --- Sheet: Confusion Matrix ---
$wsConf = $wb.Worksheets.Item("Confusion Matrix")
$wsConf.Range("B2").Value = 189
$wsConf.Range("C2").Value = 345

$wsConf.Range("B3").Value = 52
$wsConf.Range("C3").Value = 482
